# Apply the "kcal/mol -> kJ/mol" unit-conversion edit described in the
# commit message: update the column header labels and the Hartree->unit
# conversion factor used by the energy-difference formulas, tweak the H2
# cell to be a (trivial) formula, and update the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Header labels: kcal/mol -> kJ/mol (C1, E1)
$ws.Range("C1").Value = "Gas Phase Energy (kJ/mol)"
$ws.Range("E1").Value = "SMD Energy (kJ/mol)"

# 2) Conversion factor used throughout columns C and E:
#    Hartree -> kcal/mol (627.5095) becomes Hartree -> kJ/mol (2625.5)
$ws.Range("C2").Formula  = '=(B2-$B$2)*2625.5'
$ws.Range("E2").Formula  = '=(D2-$D$2)*2625.5'
$ws.Range("C3:C8").Formula = '=(B3-$B$2)*2625.5'
$ws.Range("E3:E8").Formula = '=(D3-$D$2)*2625.5'
# Re-entering the formulas can pull in the referenced cell's number format
# (B column uses a "#,##0" style); these cells were unformatted before the
# edit and stay that way, so strip any format that formula entry added.
$ws.Range("C3:C8").ClearFormats()
$ws.Range("E3:E8").ClearFormats()

# 3) H2 becomes a (trivial) formula instead of a literal constant
$ws.Range("H2").Formula = '=1.9'

# 4) Update the active selection/view
$ws.Range("H3").Select()
